$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to match the new width used by the rest of the data columns
$ws.Columns.Item(2).ColumnWidth = 15.400000000000002

# Net Income
$ws.Range("B2").Value = 1856170000.0

# Depreciation & Amortization
$ws.Range("B3").Value = 327625000.0

# Non Cash Items (Other)
$ws.Range("B4").Value = 26302000.0

# Accounts Receivable Change
$ws.Range("B5").Value = -90703000.0

# Change in inventories
$ws.Range("B6").Value = -202984000.0
$ws.Range("C6").Value = -180455000.0
$ws.Range("D6").Value = -155717000.0
$ws.Range("E6").Value = -74817000.0
$ws.Range("F6").Value = -5029000.0
$ws.Range("G6").Value = 4718000.0

# Accounts Payable Change
$ws.Range("B7").Value = 50857000.0

# Change in payables and accrued liability
$ws.Range("C8").Value = 1023000000.0
$ws.Range("D8").Value = 1018000000.0

# Change in other assets and liabilities
$ws.Range("B9").Value = 206541000.0

# Operating Cash Flow
$ws.Range("B10").Value = 2172248000.0

# Capital expenditures
$ws.Range("B11").Value = -199354000.0

# Short Term Investments Change (Net)
$ws.Range("B13").Value = 3144000.0

# Long-Term Investments Change (Net)
$ws.Range("B14").Value = -304868000.0

# Investing cash flow
$ws.Range("B16").Value = -512621000.0

# Repayment/Issuance of Debt (Net)
$ws.Range("B17").Value = -4829000.0

# Equity Repurchase (Common, Net)
$ws.Range("B18").Value = -562125000.0

# Dividends Paid (Total)
$ws.Range("B19").Value = -552848000.0

# Other financial activities
$ws.Range("B20").Value = 41710000.0

# Financing cash flow
$ws.Range("B21").Value = -1166391000.0

# Exchange Rate Adjustment
$ws.Range("B22").Value = 12739000.0

# Change in Cash
$ws.Range("B23").Value = 505975000.0

# Beginning Cash
$ws.Range("B24").Value = 1234409000.0

# Ending Cash
$ws.Range("B25").Value = 1740384000.0

# Stock Based Compensation
$ws.Range("B26").Value = 110524000.0

# Dividends Paid (Common)
$ws.Range("B27").Value = -552848000.0

# Assets Liabilities Change (Total)
$ws.Range("B28").Value = -36289000.0

# Investments Change (Net)
$ws.Range("B29").Value = -301724000.0

# Issuance/Purchase of Shares
$ws.Range("B30").Value = -562125000.0
